# Fix bug hard code in filter_date function
# Update the "last_edited_time" (column D) values for rows 2..121 on the
# single worksheet to the corrected timestamps produced after the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D13").Value   = "2024-07-04T09:47:00.000Z"
$ws.Range("D14:D26").Value  = "2024-07-04T09:48:00.000Z"
$ws.Range("D27:D46").Value  = "2024-07-04T09:47:00.000Z"
$ws.Range("D47:D80").Value  = "2024-07-04T09:48:00.000Z"
$ws.Range("D81:D120").Value = "2024-07-04T09:44:00.000Z"
$ws.Range("D121:D121").Value = "2024-07-04T09:45:00.000Z"
